# From v1.2 to v1.2.1
# Swap the step/result text between the TC2 block (row 20) and the TC3 block (row 28):
# TC2's "realizar liquidação" step moves to TC3, and TC3's "atribuir/desatribuir"
# step moves to TC2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D20").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

$ws.Range("B28").Value = "Chefe Clica para realizar a liquidação."
$ws.Range("D28").Value = "SYSTEM Apresenta a tela de Registrar Liquidações"
